$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header/data cells with their new text values.
$ws.Range("C1").Value = "Sell"
$ws.Range("D1").Value = "Ded"
$ws.Range("D2").Value = "Waget,Garnish"
$ws.Range("C2").Value = "Shirts"

# Move the active selection to C2, matching the saved view state.
$ws.Range("C2").Select()
